$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Data bestand 1"
$ws.Range("C2").Value = "dwadwaadaB"
$ws.Range("C4").Value = "dadadwadwaAB"
$ws.Range("C6").Value = "dadadadaB"
$ws.Range("C8").Value = "asddadaB"
$ws.Range("C9").Value = "wasdwaB"
$ws.Range("C10").Value = "dwadwadadaB"
$ws.Range("C11").Value = "sdwadadaB"
